$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all text-like price/percentage cells keep their exact string
# representation (leading/trailing zeros, thousand-dot separators, etc.)
# instead of being auto-coerced into numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.499.05'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.000.03'
$ws.Range('E3').Value = '  -2.92%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '537.24'
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.85'
$ws.Range('E6').Value = '  -1.33%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.994.09'
$ws.Range('E8').Value = '  -2.89%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.148'
$ws.Range('E10').Value = '  -4.99%  '
$ws.Range('E11').Value = '  -2.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.448'
$ws.Range('E12').Value = '  -2.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000222'
$ws.Range('E13').Value = '  -2.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.98'
$ws.Range('E14').Value = '  -2.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.482.88'
$ws.Range('E15').Value = '  -2.96%  '
$ws.Range('E16').Value = '  -1.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.501.56'
$ws.Range('E17').Value = '  -3.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.999.49'
$ws.Range('E18').Value = '  -2.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.63'
$ws.Range('E19').Value = '  -1.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '467.69'
$ws.Range('E20').Value = '  -4.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.22'
$ws.Range('E21').Value = '  -1.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.676'
$ws.Range('E22').Value = '  -3.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.94'
$ws.Range('E23').Value = '  -3.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.20'
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.99'
$ws.Range('E25').Value = '  -2.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('E27').Value = '  -1.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.77'
$ws.Range('E28').Value = '  -6.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('E30').Value = '  -1.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.61'
$ws.Range('E32').Value = '  +2.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.50'
$ws.Range('E33').Value = '  -0.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.29'
$ws.Range('E34').Value = '  -5.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '54.77'
$ws.Range('E35').Value = '  -3.86%  '
$ws.Range('E36').Value = '  -2.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '453.71'
$ws.Range('E37').Value = '  -8.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.168.31'
$ws.Range('E38').Value = '  -3.87%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.121'
$ws.Range('E39').Value = '  +2.59%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0787'
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0385'
$ws.Range('E41').Value = '  -3.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.14'
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.46'
$ws.Range('E43').Value = '  -7.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.99'
$ws.Range('E44').Value = '  +8.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.244'
$ws.Range('E46').Value = '  -5.51%  '
$ws.Range('E47').Value = '  -3.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '118.87'
$ws.Range('E48').Value = '  -2.35%  '
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₃0496'
$ws.Range('E50').Value = '  -7.66%  '
$ws.Range('E51').Value = '  +5.99%  '
